# issue #5: add legislator_id, name, date into dataframe
# Adds three columns - date, legislator_name, legislator_id - to the "股票"
# (stocks) sheet, mirroring the columns already present on the other sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$lastCol = 7  # existing data occupies columns A:G (total is in column G)

$dateCol = $lastCol + 1        # H
$nameCol = $lastCol + 2        # I
$idCol   = $lastCol + 3        # J

# --- Header row -----------------------------------------------------------
$ws.Cells.Item(1, $dateCol).Value = "date"
$ws.Cells.Item(1, $nameCol).Value = "legislator_name"
$ws.Cells.Item(1, $idCol).Value = "legislator_id"

# Match the look of the existing header cells (bold font + border).
$ws.Cells.Item(1, $lastCol).Copy()
$ws.Range($ws.Cells.Item(1, $dateCol), $ws.Cells.Item(1, $idCol)).PasteSpecial(-4122)

# --- Data rows --------------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    # Force the date column to stay plain text (so "2012-04-20" isn't
    # reinterpreted as a date serial number) before writing the value.
    $ws.Cells.Item($r, $dateCol).NumberFormat = "@"
    $ws.Cells.Item($r, $dateCol).Value = "2012-04-20"

    $ws.Cells.Item($r, $nameCol).Value = "蔣乃辛"
    $ws.Cells.Item($r, $idCol).Value = 1722

    # Match the look of the existing data cells for the non-date columns.
    $ws.Cells.Item($r, $lastCol).Copy()
    $ws.Range($ws.Cells.Item($r, $nameCol), $ws.Cells.Item($r, $idCol)).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
